$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A181").Value = 9999
Write-Output $ws.Range("A181").Value2
$ws.Range("M181").Value = 3
Write-Output $ws.Range("M181").Value2
